$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the simulated object (B2)
$ws.Range("B2").Value = "55Cnc_e"

# Updated input parameters (row 2)
$ws.Range("F2").Value = 0.95
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0.1
$ws.Range("L2").Value = 0.2

# Recomputed / updated output parameters (row 2)
$ws.Range("N2").Value = 0.01544
$ws.Range("O2").Value = 0.905
$ws.Range("P2").Value = 0.16192
$ws.Range("P2").NumberFormat = "0.000"
$ws.Range("Q2").Value = 0.0254
$ws.Range("R2").Value = 0.736539
$ws.Range("S2").Value = 74.3878740079582
$ws.Range("S2").NumberFormat = "0.000000000000"

# Move the AB2 helper value down to AB3 (row shifted)
$ws.Range("AB2").ClearContents()
$ws.Range("AB3").Value = 0.2828427125
$ws.Range("AB3").Style = $ws.Range("L2").Style

# Restore the active selection
$ws.Range("I3").Select()
